# Transmittals test data: include Jira RefID for each test case row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at the front of the table (shifts B:P -> C:Q)
$ws.Range("A1").EntireColumn.Insert()

# Header styling should match the other header cells (bold style, s="1")
$ws.Range("A1").Style = $ws.Range("B1").Style
$ws.Range("A1").Value = "RefID"

# Match the column width Excel computed for the new narrow RefID column
$ws.Columns("A").ColumnWidth = 10

# Jira reference IDs for each test case row
$ws.Range("A2").Value = "LATFLD-29"
$ws.Range("A3").Value = "LATFLD-35"
$ws.Range("A4").Value = "LATFLD-33"
$ws.Range("A5").Value = "LATFLD-36"

# Re-establish the (hidden) filter-database defined name scoped to this sheet,
# covering the whole table including the newly added column.
$n = $ws.Names.Add("_xlnm._FilterDatabase", "=Transmittals_Close_Cancel!`$A`$1:`$Q`$5")
$n.Visible = $false
